$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = -7.246
    "C3"  = -13.09
    "C14" = -12.509
    "C16" = -12.878
    "D18" = -8.230999999999998
    "C21" = -12.394
    "C23" = -12.446
    "D24" = -7.725999999999999
    "C25" = -12.329
    "D25" = -8.09
    "C26" = -12.947
    "D27" = -8.35
    "C29" = -12.061
    "D30" = -7.635
    "D31" = -8.251999999999999
    "D39" = -7.710000000000001
    "C40" = -12.527
    "D42" = -8.193000000000001
    "D48" = -7.100999999999999
    "D51" = -8.373999999999999
    "D52" = -7.359
    "C53" = -11.04
    "D55" = -8.145
    "D56" = -8.272
    "C57" = -13.565
    "D57" = -8.559000000000001
    "C59" = -13.077
    "D60" = -8.391999999999999
    "C65" = -12.273
    "C69" = -10.676
    "D73" = -8.004000000000001
    "D74" = -7.952999999999999
    "C79" = -12.1
    "C83" = -13.392
    "D89" = -6.858
    "D90" = -7.629
    "C91" = -10.571
    "D92" = -6.592999999999999
    "C93" = -11.786
    "C100" = -12.965
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
